# Replace the SkyWorks RF switch (U9) line item on the BOM with the
# Peregrine Semiconductor part, and update the totals section accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 45 (U9 switch): update vendor, part number, description, vendor PN and unit price ---
$ws.Range("C45").Value = "Peregrine Semiconductor"
$ws.Range("D45").Value = "PE42723A-Z"

$ws.Range("F45").Value = "UltraCMOS® SPDT RF Switch, 5–1794 MHz"
$ws.Range("F45").Style = "Normal"

$ws.Range("H45").Value = "1046-1150-1-ND"
$ws.Range("J45").Value = 4.59

# --- Totals block: push the grand-total formula down one row, leaving a blank spacer row ---
$ws.Range("L52").Formula = $ws.Range("L51").Formula
$ws.Range("L51").ClearContents()

# give the new blank J51 cell the same (currency) number format as the column above it
$ws.Range("L50").Copy()
$ws.Range("J51").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$excel.Calculate()

# restore the selection to where the edit was made
[void]$ws.Range("J46").Select()
